# The Title paragraph (first paragraph, styled "Title") is currently
# empty. Add two runs to it:
#   1. "<blockTable>" in Courier font
#   2. " Tag Demo" in the (unformatted) default run font
$d = $word.ActiveDocument
$p = $d.Paragraphs(1)
$r = $p.Range

# Remember where the paragraph's own range starts so we can address the
# freshly inserted text afterwards by absolute document position.
$start = $r.Start

# Insert both pieces of text in one go (this keeps them as a single run
# for now); we'll split the run formatting right after.
$tag = "<blockTable>"
$suffix = " Tag Demo"
$r.InsertAfter($tag + $suffix)

# Range over just "<blockTable>" -> give it the Courier font, producing
# its own <w:r><w:rPr><w:rFonts .../></w:rPr><w:t>...</w:t></w:r> run.
$rTag = $d.Range($start, $start + $tag.Length)
$rTag.Font.Name = "Courier"

# Range over " Tag Demo" -> touch (and revert) a character property so
# Word materializes an explicit, but empty, <w:rPr/> on its own run
# instead of merging back into the Courier run.
$rSuffix = $d.Range($start + $tag.Length, $start + $tag.Length + $suffix.Length)
$rSuffix.Bold = 1
$rSuffix.Bold = 0
